# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
# Both sheets carry the same event rows, so the same F-column updates apply
# to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3053
    7  = 1670
    12 = 1373
    14 = 515
    16 = 32
    23 = 3203
    25 = 134
    26 = 315
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
